# Apply the "Add files via upload" update to solieu.xlsx
#  - Sheet "DANH SÁCH NỢ": row 6 (A7=6) is now fully paid -> mark it
#    "Đã trả đủ" and highlight the row green; add three new debt rows
#    (A10=9, A11=10, A12=11).
#  - Sheet "THONG KE NAP ": log the matching payment entries in rows
#    106-109.
#  - Selections move to reflect where the user was last working.

$wb = $excel.ActiveWorkbook

$wsDanhSach = $wb.Worksheets.Item("DANH SÁCH NỢ")
$wsThongKeNap = $wb.Worksheets.Item("THONG KE NAP ")

# ---------------------------------------------------------------------
# 1) Row 7 (STT 6) gets paid off in full: the "phải trả"(E) column now
#    shows the 20000 payment, the due date (K) moves up, and the status
#    (M) flips from "Chưa trả đủ" to "Đã trả đủ".
# ---------------------------------------------------------------------
$wsDanhSach.Range("E7").Value = 20000
$wsDanhSach.Range("K7").Value = 46016
$wsDanhSach.Range("M7").Value = "Đã trả đủ"

# Highlight the whole paid-off row with a green fill.
$wsDanhSach.Range("B7:M7").Interior.Color = 5296274

# ---------------------------------------------------------------------
# 2) Three brand-new debt entries fill rows 10-12 (previously blank).
# ---------------------------------------------------------------------

# Row 10 - STT 9
$wsDanhSach.Range("B10").Value = "Bùi Anh Tài"
$wsDanhSach.Range("C10").Value = "Nạp Robux"
$wsDanhSach.Range("D10").Value = 50000
$wsDanhSach.Range("E10").Value = 0
$wsDanhSach.Range("F10").Formula = "=(D10+I10)-E10"
$wsDanhSach.Range("G10").Value = 0
$wsDanhSach.Range("H10").Value = 0
$wsDanhSach.Range("I10").Formula = "=D10*H10"
$wsDanhSach.Range("J10").Value = 46016
$wsDanhSach.Range("K10").Value = 46022
$wsDanhSach.Range("M10").Value = "Chưa trả đủ"

# Row 11 - STT 10
$wsDanhSach.Range("B11").Value = "Trần Huỳnh Như Ý"
$wsDanhSach.Range("C11").Value = "Mua thẻ Zing"
$wsDanhSach.Range("D11").Value = 20000
$wsDanhSach.Range("E11").Value = 0
$wsDanhSach.Range("F11").Formula = "=(D11+I11)-E11"
$wsDanhSach.Range("G11").Value = 0
$wsDanhSach.Range("H11").Value = 0
$wsDanhSach.Range("I11").Formula = "=D11*H11"
$wsDanhSach.Range("J11").Value = 46017
$wsDanhSach.Range("K11").Value = 46023
$wsDanhSach.Range("M11").Value = "Chưa trả đủ"

# Row 12 - STT 11
$wsDanhSach.Range("B12").Value = "Bùi Bích Ngọc"
$wsDanhSach.Range("C12").Value = "Mua thẻ Zing"
$wsDanhSach.Range("D12").Value = 20000
$wsDanhSach.Range("E12").Value = 0
$wsDanhSach.Range("F12").Formula = "=(D12+I12)-E12"
$wsDanhSach.Range("G12").Value = 0
$wsDanhSach.Range("H12").Value = 0
$wsDanhSach.Range("I12").Formula = "=D12*H12"
$wsDanhSach.Range("J12").Value = 46017
$wsDanhSach.Range("K12").Value = 46023
$wsDanhSach.Range("M12").Value = "Chưa trả đủ"

# ---------------------------------------------------------------------
# 3) Mirror the new payments into the "THONG KE NAP " log (rows 106-109).
# ---------------------------------------------------------------------
$wsThongKeNap.Range("A106").Value = 46016
$wsThongKeNap.Range("B106").Value = "Nguyễn Huy Hoàng"
$wsThongKeNap.Range("C106").Value = 20000
$wsThongKeNap.Range("D106").Value = "Nạp Free Fire"

$wsThongKeNap.Range("A107").Value = 46016
$wsThongKeNap.Range("B107").Value = "Bùi Anh Tài"
$wsThongKeNap.Range("C107").Value = 50000
$wsThongKeNap.Range("D107").Value = "Nạp Robux"

$wsThongKeNap.Range("A108").Value = 46017
$wsThongKeNap.Range("B108").Value = "Trần Huỳnh Như Ý"
$wsThongKeNap.Range("C108").Value = 20000
$wsThongKeNap.Range("D108").Value = "Mua thẻ Zing"

$wsThongKeNap.Range("A109").Value = 46017
$wsThongKeNap.Range("B109").Value = "Bùi Bích Ngọc"
$wsThongKeNap.Range("C109").Value = 20000
$wsThongKeNap.Range("D109").Value = "Mua thẻ Zing"

# ---------------------------------------------------------------------
# 4) Update the cursor/selection on each sheet, ending back on the
#    originally-active "DANH SÁCH NỢ" tab.
# ---------------------------------------------------------------------
$wsThongKeNap.Activate()
$wsThongKeNap.Range("D110").Select()

$wsDanhSach.Activate()
$wsDanhSach.Range("I24").Select()
